$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(12, 8).Value = 356.2857  # H12
$ws.Cells.Item(12, 9).Value = 338.8  # I12
$ws.Cells.Item(12, 11).Value = 338.8  # K12
$ws.Cells.Item(12, 13).Value = -168.8  # M12
$ws.Cells.Item(18, 8).Value = 1490.6923  # H18
$ws.Cells.Item(18, 10).Value = 1844.75  # J18
$ws.Cells.Item(18, 12).Value = 1844.75  # L18
$ws.Cells.Item(18, 14).Value = -2412.75  # N18
$ws.Cells.Item(19, 8).Value = 583.5  # H19
$ws.Cells.Item(19, 9).Value = 767  # I19
$ws.Cells.Item(19, 10).Value = 400  # J19
$ws.Cells.Item(19, 11).Value = 767  # K19
$ws.Cells.Item(19, 12).Value = 400  # L19
$ws.Cells.Item(19, 13).Value = -592  # M19
$ws.Cells.Item(19, 14).Value = -750  # N19
$ws.Cells.Item(33, 8).Value = 462.33334  # H33
$ws.Cells.Item(33, 9).Value = 431.72223  # I33
$ws.Cells.Item(33, 10).Value = 554.1667  # J33
$ws.Cells.Item(33, 11).Value = 431.72223  # K33
$ws.Cells.Item(33, 12).Value = 554.1667  # L33
$ws.Cells.Item(33, 13).Value = -202.72223  # M33
$ws.Cells.Item(33, 14).Value = -1012.1667  # N33
$ws.Cells.Item(43, 8).Value = 5567195.5  # H43
$ws.Cells.Item(43, 9).Value = 20720.2  # I43
$ws.Cells.Item(43, 10).Value = 11113671  # J43
$ws.Cells.Item(43, 11).Value = 20720.2  # K43
$ws.Cells.Item(43, 12).Value = 11113671  # L43
$ws.Cells.Item(43, 13).Value = -20651.2  # M43
$ws.Cells.Item(43, 14).Value = -11113809  # N43
$ws.Cells.Item(62, 8).Value = 13892511  # H62
$ws.Cells.Item(62, 9).Value = 15876442  # I62
$ws.Cells.Item(62, 11).Value = 15876442  # K62
$ws.Cells.Item(62, 13).Value = -15875818  # M62
$ws.Cells.Item(65, 8).Value = 13892511  # H65
$ws.Cells.Item(65, 9).Value = 15876442  # I65
$ws.Cells.Item(65, 11).Value = 79382210  # K65
$ws.Cells.Item(65, 13).Value = -79379090  # M65
$ws.Cells.Item(111, 8).Value = 1980  # H111
$ws.Cells.Item(111, 9).Value = 1980  # I111
$ws.Cells.Item(111, 10).Value = 0  # J111
$ws.Cells.Item(111, 11).Value = 5940  # K111
$ws.Cells.Item(111, 12).Value = 0  # L111
$ws.Cells.Item(111, 13).Value = -2873  # M111
$ws.Cells.Item(111, 14).ClearContents()  # N111 (was -25634)
$ws.Cells.Item(113, 8).Value = 3364.3333  # H113
$ws.Cells.Item(113, 9).Value = 3493.3333  # I113
$ws.Cells.Item(113, 11).Value = 3493.3333  # K113
$ws.Cells.Item(113, 13).Value = -239.3332999999998  # M113
$ws.Cells.Item(116, 8).Value = 3360.0625  # H116
$ws.Cells.Item(116, 10).Value = 3955.1428  # J116
$ws.Cells.Item(116, 12).Value = 3955.1428  # L116
$ws.Cells.Item(116, 14).Value = -10839.1428  # N116
$ws.Cells.Item(131, 8).Value = 880  # H131
$ws.Cells.Item(131, 9).Value = 880  # I131
$ws.Cells.Item(131, 11).Value = 2640  # K131
$ws.Cells.Item(131, 13).Value = 2400  # M131
$ws.Cells.Item(137, 8).Value = 1566.5294  # H137
$ws.Cells.Item(137, 9).Value = 1160.2354  # I137
$ws.Cells.Item(137, 10).Value = 1972.8235  # J137
$ws.Cells.Item(137, 11).Value = 3480.7062  # K137
$ws.Cells.Item(137, 12).Value = 5918.470499999999  # L137
$ws.Cells.Item(137, 13).Value = -930.7062000000001  # M137
$ws.Cells.Item(137, 14).Value = -11018.4705  # N137

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 4701.3945  # H32
$ws.Cells.Item(32, 9).Value = 4526.4707  # I32
$ws.Cells.Item(32, 11).Value = 4526.4707  # K32
$ws.Cells.Item(32, 13).Value = -4239.4707  # M32
$ws.Cells.Item(45, 8).Value = 1153.8  # H45
$ws.Cells.Item(45, 9).Value = 1089.1428  # I45
$ws.Cells.Item(45, 10).Value = 1304.6666  # J45
$ws.Cells.Item(45, 11).Value = 1089.1428  # K45
$ws.Cells.Item(45, 12).Value = 1304.6666  # L45
$ws.Cells.Item(45, 13).Value = -712.1428000000001  # M45
$ws.Cells.Item(45, 14).Value = -2058.6666  # N45
$ws.Cells.Item(74, 8).Value = 2710.0908  # H74
$ws.Cells.Item(74, 9).Value = 1726.375  # I74
$ws.Cells.Item(74, 11).Value = 1726.375  # K74
$ws.Cells.Item(74, 13).Value = -852.375  # M74
$ws.Cells.Item(77, 8).Value = 2710.0908  # H77
$ws.Cells.Item(77, 9).Value = 1726.375  # I77
$ws.Cells.Item(77, 11).Value = 8631.875  # K77
$ws.Cells.Item(77, 13).Value = -4263.875  # M77
$ws.Cells.Item(122, 8).Value = 2215.7693  # H122
$ws.Cells.Item(122, 10).Value = 2452.1428  # J122
$ws.Cells.Item(122, 12).Value = 7356.428400000001  # L122
$ws.Cells.Item(122, 14).Value = -12256.4284  # N122

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(25, 8).Value = 580  # H25
$ws.Cells.Item(25, 9).Value = 580  # I25
$ws.Cells.Item(25, 11).Value = 580  # K25
$ws.Cells.Item(25, 13).Value = -345  # M25
$ws.Cells.Item(87, 8).Value = 53000  # H87
$ws.Cells.Item(87, 10).Value = 53000  # J87
$ws.Cells.Item(87, 12).Value = 53000  # L87
$ws.Cells.Item(87, 14).Value = -55496  # N87
$ws.Cells.Item(90, 8).Value = 53000  # H90
$ws.Cells.Item(90, 10).Value = 53000  # J90
$ws.Cells.Item(90, 12).Value = 159000  # L90
$ws.Cells.Item(90, 14).Value = -171480  # N90
$ws.Cells.Item(139, 8).Value = 36353.332  # H139
$ws.Cells.Item(139, 10).Value = 36353.332  # J139
$ws.Cells.Item(139, 12).Value = 36353.332  # L139
$ws.Cells.Item(139, 14).Value = -46633.332  # N139

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 45455630  # H16
$ws.Cells.Item(16, 9).Value = 66667770  # I16
$ws.Cells.Item(16, 11).Value = 66667770  # K16
$ws.Cells.Item(16, 13).Value = -66667483  # M16
$ws.Cells.Item(31, 8).Value = 1894.3077  # H31
$ws.Cells.Item(31, 9).Value = 1894.3077  # I31
$ws.Cells.Item(31, 11).Value = 1894.3077  # K31
$ws.Cells.Item(31, 13).Value = -1599.3077  # M31
$ws.Cells.Item(32, 8).Value = 0  # H32
$ws.Cells.Item(32, 9).Value = 0  # I32
$ws.Cells.Item(32, 11).Value = 0  # K32
$ws.Cells.Item(32, 13).ClearContents()  # M32 (was -4684)
$ws.Cells.Item(34, 8).Value = 1894.3077  # H34
$ws.Cells.Item(34, 9).Value = 1894.3077  # I34
$ws.Cells.Item(34, 11).Value = 1894.3077  # K34
$ws.Cells.Item(34, 13).Value = -1692.3077  # M34
$ws.Cells.Item(74, 8).Value = 32500  # H74
$ws.Cells.Item(74, 10).Value = 32500  # J74
$ws.Cells.Item(74, 12).Value = 32500  # L74
$ws.Cells.Item(74, 14).Value = -34248  # N74
$ws.Cells.Item(77, 8).Value = 32500  # H77
$ws.Cells.Item(77, 10).Value = 32500  # J77
$ws.Cells.Item(77, 12).Value = 97500  # L77
$ws.Cells.Item(77, 14).Value = -106236  # N77
$ws.Cells.Item(86, 8).Value = 6112515.5  # H86
$ws.Cells.Item(86, 9).Value = 9561137  # I86
$ws.Cells.Item(86, 11).Value = 9561137  # K86
$ws.Cells.Item(86, 13).Value = -9560014  # M86
$ws.Cells.Item(89, 8).Value = 6112515.5  # H89
$ws.Cells.Item(89, 9).Value = 9561137  # I89
$ws.Cells.Item(89, 11).Value = 47805685  # K89
$ws.Cells.Item(89, 13).Value = -47800069  # M89
$ws.Cells.Item(108, 8).Value = 32401  # H108
$ws.Cells.Item(108, 10).Value = 32401  # J108
$ws.Cells.Item(108, 12).Value = 32401  # L108
$ws.Cells.Item(108, 14).Value = -40081  # N108
$ws.Cells.Item(113, 8).Value = 45455630  # H113
$ws.Cells.Item(113, 9).Value = 66667770  # I113
$ws.Cells.Item(113, 11).Value = 66667770  # K113
$ws.Cells.Item(113, 13).Value = -66665600  # M113
$ws.Cells.Item(132, 8).Value = 3230.75  # H132
$ws.Cells.Item(132, 9).Value = 3061.8333  # I132
$ws.Cells.Item(132, 10).Value = 3399.6667  # J132
$ws.Cells.Item(132, 11).Value = 9185.499899999999  # K132
$ws.Cells.Item(132, 12).Value = 10199.0001  # L132
$ws.Cells.Item(132, 13).Value = -6655.499899999999  # M132
$ws.Cells.Item(132, 14).Value = -15259.0001  # N132
$ws.Cells.Item(141, 8).Value = 765152.5  # H141
$ws.Cells.Item(141, 10).Value = 765152.5  # J141
$ws.Cells.Item(141, 12).Value = 765152.5  # L141
$ws.Cells.Item(141, 14).Value = -775512.5  # N141

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(70, 8).Value = 4769.1177  # H70
$ws.Cells.Item(70, 10).Value = 5966.5835  # J70
$ws.Cells.Item(70, 12).Value = 17899.7505  # L70
$ws.Cells.Item(70, 14).Value = -18529.7505  # N70
$ws.Cells.Item(73, 8).Value = 4769.1177  # H73
$ws.Cells.Item(73, 10).Value = 5966.5835  # J73
$ws.Cells.Item(73, 12).Value = 17899.7505  # L73
$ws.Cells.Item(73, 14).Value = -20083.7505  # N73
$ws.Cells.Item(131, 8).Value = 10527373  # H131
$ws.Cells.Item(131, 10).Value = 1102.6405  # J131
$ws.Cells.Item(131, 12).Value = 3307.9215  # L131
$ws.Cells.Item(131, 14).Value = -13387.9215  # N131

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(107, 8).Value = 612.2727  # H107
$ws.Cells.Item(107, 9).Value = 556.4  # I107
$ws.Cells.Item(107, 10).Value = 658.8333  # J107
$ws.Cells.Item(107, 11).Value = 556.4  # K107
$ws.Cells.Item(107, 12).Value = 658.8333  # L107
$ws.Cells.Item(107, 13).Value = 1363.6  # M107
$ws.Cells.Item(107, 14).Value = -4498.8333  # N107
$ws.Cells.Item(113, 8).Value = 1404.625  # H113
$ws.Cells.Item(113, 9).Value = 1164.5555  # I113
$ws.Cells.Item(113, 11).Value = 1164.5555  # K113
$ws.Cells.Item(113, 13).Value = 1005.4445  # M113

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 592.2  # H16
$ws.Cells.Item(16, 9).Value = 607.5789  # I16
$ws.Cells.Item(16, 10).Value = 300  # J16
$ws.Cells.Item(16, 11).Value = 607.5789  # K16
$ws.Cells.Item(16, 12).Value = 300  # L16
$ws.Cells.Item(16, 13).Value = -437.5789  # M16
$ws.Cells.Item(16, 14).Value = -640  # N16
$ws.Cells.Item(132, 8).Value = 3255.2  # H132
$ws.Cells.Item(132, 9).Value = 3270.6  # I132
$ws.Cells.Item(132, 10).Value = 3239.8  # J132
$ws.Cells.Item(132, 11).Value = 9811.799999999999  # K132
$ws.Cells.Item(132, 12).Value = 9719.400000000001  # L132
$ws.Cells.Item(132, 13).Value = -7281.799999999999  # M132
$ws.Cells.Item(132, 14).Value = -14779.4  # N132

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 1400  # H81
$ws.Cells.Item(84, 8).Value = 1400  # H84
$ws.Cells.Item(132, 8).Value = 1212.8206  # H132
$ws.Cells.Item(132, 9).Value = 1059.0322  # I132
$ws.Cells.Item(132, 10).Value = 1808.75  # J132
$ws.Cells.Item(132, 11).Value = 3177.0966  # K132
$ws.Cells.Item(132, 12).Value = 5426.25  # L132
$ws.Cells.Item(132, 13).Value = -647.0966000000003  # M132
$ws.Cells.Item(132, 14).Value = -10486.25  # N132
